$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price"/"Volume(1h)" figures published by the latest
# coinranking.com snapshot (including the TRON / WrappedEther rank swap).

# Row 2
$ws.Range("D2").Value = '27.346.34'
$ws.Range("E2").Value = '  +1.30%  '
# Row 3
$ws.Range("D3").Value = '1.857.47'
$ws.Range("E3").Value = '  +1.60%  '
# Row 4
$ws.Range("E4").Value = '  -0.79%  '
# Row 5
$ws.Range("D5").Value = '''314.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.73%  '
# Row 7
$ws.Range("D7").Value = '''0.4615'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.75%  '
# Row 8
$ws.Range("D8").Value = '''0.3713'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.31%  '
# Row 9
$ws.Range("D9").Value = '''0.07314'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.61%  '
# Row 10
$ws.Range("D10").Value = '''0.8807'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.78%  '
# Row 11
$ws.Range("D11").Value = '''19.89'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '
# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.957.29'
$ws.Range("E12").Value = '  +5.01%  '
# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.07790'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.23%  '
# Row 14
$ws.Range("D14").Value = '''5.387'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.57%  '
# Row 15
$ws.Range("D15").Value = '''6.549'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.47%  '
# Row 16
$ws.Range("D16").Value = '''91.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '
# Row 17
$ws.Range("D17").Value = '''1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.94%  '
# Row 18
$ws.Range("D18").Value = '''0.000009095'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.50%  '
# Row 19
$ws.Range("E19").Value = '  -0.68%  '
# Row 20
$ws.Range("E20").Value = '  +0.62%  '
# Row 21
$ws.Range("D21").Value = '27.353.71'
$ws.Range("E21").Value = '  +1.61%  '
# Row 22
$ws.Range("D22").Value = '''5.129'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.53%  '
# Row 23
$ws.Range("D23").Value = '''10.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.60%  '
# Row 24
$ws.Range("D24").Value = '2.140.13'
$ws.Range("E24").Value = '  +2.54%  '
# Row 25
$ws.Range("D25").Value = '''1.930'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.45%  '
# Row 26
$ws.Range("D26").Value = '''152.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.30%  '
# Row 27
$ws.Range("D27").Value = '''18.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.51%  '
# Row 28
$ws.Range("D28").Value = '''2.071'
$ws.Range("D28").Style = "Normal"
# Row 29
$ws.Range("D29").Value = '''5.104'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.43%  '
# Row 30
$ws.Range("D30").Value = '''116.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.51%  '
# Row 31
$ws.Range("D31").Value = '''0.08859'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '
# Row 32
$ws.Range("D32").Value = '''0.7726'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.19%  '
# Row 33
$ws.Range("D33").Value = '''3.039'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.89%  '
# Row 35
$ws.Range("D35").Value = '''4.494'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.15%  '
# Row 36
$ws.Range("D36").Value = '''2.657'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.00%  '
# Row 37
$ws.Range("E37").Value = '  +0.12%  '
# Row 39
$ws.Range("D39").Value = '''0.05227'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.16%  '
# Row 40
$ws.Range("E40").Value = '  +0.78%  '
# Row 41
$ws.Range("D41").Value = '''7.014'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.99%  '
# Row 42
$ws.Range("D42").Value = '''0.5144'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.79%  '
# Row 43
$ws.Range("D43").Value = '''0.1636'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.72%  '
# Row 44
$ws.Range("D44").Value = '''8.404'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.37%  '
# Row 45
$ws.Range("D45").Value = '''0.4827'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.23%  '
# Row 46
$ws.Range("D46").Value = '''10.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.87%  '
# Row 47
$ws.Range("E47").Value = '  -0.80%  '
# Row 48
$ws.Range("D48").Value = '''102.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.17%  '
# Row 49
$ws.Range("D49").Value = '''1.652'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.75%  '
# Row 50
$ws.Range("D50").Value = '''0.06218'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.08%  '
# Row 51
$ws.Range("D51").Value = '''65.65'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.10%  '
